$wb = $excel.ActiveWorkbook

# Sheet 1: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L18").Value = 4413.66

# Sheet 2: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F18").Value = 4413.66
$ws2.Range("F29").Value = 7536.18

# Sheet 3: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D15").Value = 6084.09
$ws3.Range("E15").Value = -4446.09
$ws3.Range("F15").Value = 3.714340659340659
$ws3.Range("D19").Value = 7536.179999999999
$ws3.Range("E19").Value = 29963.82093005039
$ws3.Range("F19").Value = 0.200964795015803
